$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("398898064313548820", "3", "11", "Leboge#3392"),
    @("338028125269000199", "8", "9", "Reptile#9182"),
    @("552891291357282304", "8", "8", "Cyber#5379")
)

$r = 32
foreach ($row in $data) {
    for ($col = 1; $col -le 4; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $value = $row[$col - 1]
        if ($value -match '^[0-9]+$') {
            # Force text storage for numeric-looking values (IDs, month, day)
            # without leaving a visible/applied style on the cell.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
    $r++
}
